# Updated cryptos list on Thu Mar 30 21:21:03 UTC 2023 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $val) {
    $c = $ws.Range($addr)
    $c.NumberFormat = "@"
    $c.Value2 = $val
}

# Row 2 - Bitcoin
Set-TextValue "D2" "28.132.05"
$ws.Range("E2").Value2 = "  -1.00%  "

# Row 3 - Ethereum
Set-TextValue "D3" "1.794.47"
$ws.Range("E3").Value2 = "  -0.50%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value2 = "  +0.01%  "

# Row 5 - BNB
Set-TextValue "D5" "317.03"
$ws.Range("E5").Value2 = "  +0.51%  "

# Row 6 - USDC
Set-TextValue "D6" "0.9999"
$ws.Range("E6").Value2 = "  -0.05%  "

# Row 7 - XRP
Set-TextValue "D7" "0.5355"
$ws.Range("E7").Value2 = "  -2.35%  "

# Row 8 - Cardano
Set-TextValue "D8" "0.3765"
$ws.Range("E8").Value2 = "  -2.39%  "

# Row 9 - Dogecoin
Set-TextValue "D9" "0.07474"
$ws.Range("E9").Value2 = "  -1.90%  "

# Row 10 - OKB
Set-TextValue "D10" "41.74"
$ws.Range("E10").Value2 = "  -1.86%  "

# Row 11 - Polygon
$ws.Range("E11").Value2 = "  -2.85%  "

# Row 13 - Solana
$ws.Range("E13").Value2 = "  -2.99%  "

# Row 14 - Polkadot
Set-TextValue "D14" "6.107"
$ws.Range("E14").Value2 = "  -1.28%  "

# Rows 15 & 16 swap: Chainlink moves to row 15, WrappedEther moves to row 16
$ws.Range("B15").Value2 = "Chainlink"
$ws.Range("C15").Value2 = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
Set-TextValue "D15" "7.220"
$ws.Range("E15").Value2 = "  -2.91%  "

$ws.Range("B16").Value2 = "WrappedEther"
$ws.Range("C16").Value2 = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
Set-TextValue "D16" "1.783.80"
$ws.Range("E16").Value2 = "  -1.12%  "

# Row 17 - Litecoin
Set-TextValue "D17" "89.17"
$ws.Range("E17").Value2 = "  -3.18%  "

# Row 18 - ShibaInu
$ws.Range("E18").Value2 = "  -1.30%  "

# Row 19 - TRON
Set-TextValue "D19" "0.06454"
$ws.Range("E19").Value2 = "  +0.31%  "

# Row 20 - Dai
Set-TextValue "D20" "0.9994"
$ws.Range("E20").Value2 = "  -0.08%  "

# Row 21 - Avalanche
Set-TextValue "D21" "17.35"
$ws.Range("E21").Value2 = "  +0.14%  "

# Row 22 - Uniswap
Set-TextValue "D22" "5.904"
$ws.Range("E22").Value2 = "  -1.16%  "

# Row 23 - WrappedBTC
Set-TextValue "D23" "28.155.80"
$ws.Range("E23").Value2 = "  -0.93%  "

# Row 24 - Cosmos
Set-TextValue "D24" "11.19"
$ws.Range("E24").Value2 = "  -2.04%  "

# Row 25 - Toncoin
Set-TextValue "D25" "2.093"
$ws.Range("E25").Value2 = "  -2.07%  "

# Row 26 - Monero
Set-TextValue "D26" "154.83"
$ws.Range("E26").Value2 = "  -2.60%  "

# Row 27 - EthereumClassic
Set-TextValue "D27" "20.22"
$ws.Range("E27").Value2 = "  -2.32%  "

# Row 28 - WrappedliquidstakedEther2.0
Set-TextValue "D28" "1.990.72"
$ws.Range("E28").Value2 = "  -1.12%  "

# Row 29 - LidoDAOToken
Set-TextValue "D29" "2.282"
$ws.Range("E29").Value2 = "  -5.28%  "

# Row 30 - BitcoinCash
Set-TextValue "D30" "120.58"
$ws.Range("E30").Value2 = "  -2.74%  "

# Row 31 - ImmutableX
$ws.Range("E31").Value2 = "  -0.61%  "

# Row 32 - Stellar
Set-TextValue "D32" "0.1051"
$ws.Range("E32").Value2 = "  +2.95%  "

# Row 33 - HuobiToken
Set-TextValue "D33" "3.652"
$ws.Range("E33").Value2 = "  -0.94%  "

# Row 34 - Filecoin
Set-TextValue "D34" "5.565"
$ws.Range("E34").Value2 = "  -3.83%  "

# Row 35 - Hedera
Set-TextValue "D35" "0.06535"
$ws.Range("E35").Value2 = "  +1.34%  "

# Row 37 - VeChain
Set-TextValue "D37" "0.02284"

# Row 38 - InternetComputer(DFINITY)
Set-TextValue "D38" "5.038"
$ws.Range("E38").Value2 = "  -2.89%  "

# Row 39 - FraxShare
$ws.Range("E39").Value2 = "  -4.12%  "

# Row 40 - WEMIXTOKEN
Set-TextValue "D40" "1.449"
$ws.Range("E40").Value2 = "  +4.51%  "

# Row 41 - TheSandbox
Set-TextValue "D41" "0.6166"
$ws.Range("E41").Value2 = "  -3.93%  "

# Row 42 - Aptos
Set-TextValue "D42" "11.09"
$ws.Range("E42").Value2 = "  -5.01%  "

# Row 43 - TrustWalletToken
Set-TextValue "D43" "1.173"
$ws.Range("E43").Value2 = "  +0.91%  "

# Row 44 - Frax
Set-TextValue "D44" "0.9992"
$ws.Range("E44").Value2 = "  -0.09%  "

# Row 45 - EnergySwap
Set-TextValue "D45" "13.29"
$ws.Range("E45").Value2 = "  -2.33%  "

# Row 46 - PancakeSwap
$ws.Range("E46").Value2 = "  -0.13%  "

# Row 47 - Decentraland
Set-TextValue "D47" "0.5784"
$ws.Range("E47").Value2 = "  -3.37%  "

# Row 48 - Quant
Set-TextValue "D48" "127.34"
$ws.Range("E48").Value2 = "  +0.27%  "

# Row 49 - EOS
Set-TextValue "D49" "1.190"
$ws.Range("E49").Value2 = "  +3.35%  "

# Row 50 - NEARProtocol
$ws.Range("E50").Value2 = "  -2.78%  "

# Row 51 - Cronos
Set-TextValue "D51" "0.06815"
$ws.Range("E51").Value2 = "  -1.11%  "
